$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.098.95"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "1.650.20"
$ws.Range("E3").Value = "  -4.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4792"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.71%  "
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06033"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07092"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "1.652.26"
$ws.Range("E11").Value = "  -4.89%  "
$ws.Range("E12").Value = "  -3.36%  "
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.561"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "25.079.43"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006547"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.406"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").Value = "1.861.23"
$ws.Range("E22").Value = "  -5.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.464"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.234"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "133.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.396"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.691"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "101.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.787"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04538"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9421"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5784"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.620"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.0000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3702"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.787"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1128"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.025"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05158"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "51.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3335"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.27%  "
